# Updates odds/values in the FlashScore weekly games workbook to match
# the new scraped data snapshot. Only numeric value cells change; no
# structural, formatting, or text changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "I3"  = 2.9
    "Q3"  = 2.6
    "R3"  = 1.48
    "Z3"  = 26
    "AP3" = 29
    "AV3" = 67
    "AZ3" = 51

    "M4"  = 1.13
    "N4"  = 6
    "O4"  = 1.53
    "P4"  = 2.38
    "U4"  = 2.38
    "V4"  = 1.53
    "Y4"  = 9.5
    "AC4" = 6
    "AG4" = 9.5

    "Q6" = 1.9
    "R6" = 1.95

    "H7"  = 2.95
    "I7"  = 3.35
    "J7"  = 2.82
    "K7"  = 1.98
    "L7"  = 3.8
    "P7"  = 2.6
    "T7"  = 2.45
    "U7"  = 1.82
    "W7"  = 6.4
    "X7"  = 9.75
    "Y7"  = 9
    "AA7" = 20
    "AB7" = 35
    "AC7" = 7.6
    "AD7" = 5.8
    "AG7" = 9.25
    "AH7" = 18
    "AI7" = 11.5
    "AK7" = 32
    "AL7" = 40
    "AN7" = 4
    "AO7" = 11.75
    "AP7" = 21
    "AQ7" = 50
    "AR7" = 90
    "AT7" = 2.42
    "AU7" = 6.8
    "AX7" = 18
    "AY7" = 24
    "AZ7" = 90

    "M8" = 1.07
    "N8" = 9

    "G13"  = 1.75
    "H13"  = 3.25
    "I13"  = 4.75
    "T13"  = 2.4
    "W13"  = 5.6
    "X13"  = 7.3
    "AA13" = 15.5
    "AC13" = 7.6
    "AD13" = 6.4
    "AF13" = 110
    "AG13" = 10.75
    "AH13" = 26
    "AL13" = 65
    "AN13" = 3.45
    "AP13" = 19.5
    "AT13" = 2.37
    "BA13" = 250
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
